$d = $word.ActiveDocument

# The paragraph contains: Service(service_id, service_description, tag, url, script, callback)
# The run ", service_description, tag, url, script, callback" needs to be split into three runs:
#   1. ", service_description, tag, url, "   (keeps trailing space -> xml:space="preserve")
#   2. "script"
#   3. ", callback"
# The following ")" run must stay untouched.

$wholeText = ", service_description, tag, url, script, callback"
$seg1 = ", service_description, tag, url, "
$seg2 = "script"
$seg3 = ", callback"

$find = $d.Content.Find
$found = $find.Execute($wholeText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $whole = $find.Parent.Duplicate

    $start = $whole.Start
    $b1 = $start + $seg1.Length
    $b2 = $b1 + $seg2.Length
    $end = $b2 + $seg3.Length

    $r1 = $d.Range($start, $b1)
    $r2 = $d.Range($b1, $b2)
    $r3 = $d.Range($b2, $end)

    # Re-assigning a range's own FormattedText forces the run to become its
    # own independent run (splitting it off from neighboring text) without
    # altering any character formatting.
    $r1.FormattedText = $r1.FormattedText
    $r2.FormattedText = $r2.FormattedText
    $r3.FormattedText = $r3.FormattedText
}
